# Replace use of escape() with encodeURIComponent() in the "queries" sheet's
# auxillaryHash column (G2:G4).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("queries")

$newFormula = "'household_id='+encodeURIComponent(data('household_id'))"

$ws.Range("G2").Value = $newFormula
$ws.Range("G3").Value = $newFormula
$ws.Range("G4").Value = $newFormula

# Restore the originally active sheet/selection (household sheet with H7 selected)
$wsHousehold = $wb.Worksheets.Item("household")
$wsHousehold.Activate()
$wsHousehold.Range("H7").Select()
